# "impl advantage filter for excel data"
# Move the ID/First Name/Last Name/Country table from E13:H23 up to A1:D11
# so it starts at the top-left of the sheet (ready for a data filter/table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$src = $ws.Range("E13:H23")
$dst = $ws.Range("A1:D11")

# Cut the table (values + formatting, incl. the bold header style) and
# drop it at its new home.
$src.Cut($dst)

# Clear out whatever is left behind at the old location.
$ws.Range("E13:H23").Clear()

# Select the relocated table, matching the new selection state.
$ws.Range("A1:D11").Select()
